$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns at S:T, shifting the existing piezo_y..c block
# (previously S:X) two columns to the right (becomes U:Z).
$ws.Range("S1:T1").EntireColumn.Insert()

# New column headers for the inserted columns.
$ws.Range("S1").Value = "shear_reinf"
$ws.Range("T1").Value = "normal_reinf"

# Fill the new columns with 0 for every data row (2-21).
$ws.Range("S2:T21").Value = 0
